# Quarterly income-statement update: shift every period one column to the
# left (drop the oldest "6 ماهه منتهی به 1399/06" period / its publish
# date) and append the new "12 ماهه منتهی به 1401/12" period (published
# 1402-02-28) in column M. Also corrects a handful of figures produced by
# the revised read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 8: financial-period headers (D8:M8) - shift left, new period in M8
# ---------------------------------------------------------------------
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------
# Row 9: publish-date headers (D9:M9) - shift left, new date in M9
# ---------------------------------------------------------------------
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-03-25 (8)"
$ws.Range("F9").Value = "1401-04-28 (2)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-28 (7)"
$ws.Range("J9").Value = "1401-04-28"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-28"

# ---------------------------------------------------------------------
# Row 11: فروش (Sales)
# ---------------------------------------------------------------------
$ws.Range("D11").Value = 5144284
$ws.Range("E11").Value = 7130495
$ws.Range("F11").Value = 2669782
$ws.Range("G11").Value = 5639882
$ws.Range("H11").Value = 8256199
$ws.Range("I11").Value = 10532870
$ws.Range("J11").Value = 3911403
$ws.Range("K11").Value = 8146029
$ws.Range("L11").Value = 11997322
$ws.Range("M11").Value = 17595389

# ---------------------------------------------------------------------
# Row 12: بهای تمام شده کالای فروش رفته (COGS)
# ---------------------------------------------------------------------
$ws.Range("D12").Value = -3115427
$ws.Range("E12").Value = -4662621
$ws.Range("F12").Value = -2087128
$ws.Range("G12").Value = -4447095
$ws.Range("H12").Value = -6562693
$ws.Range("I12").Value = -8613656
$ws.Range("J12").Value = -3012110
$ws.Range("K12").Value = -6415675
$ws.Range("L12").Value = -9460322
$ws.Range("M12").Value = -13991314

# ---------------------------------------------------------------------
# Row 13: سود (زیان) ناخالص (Gross profit)
# ---------------------------------------------------------------------
$ws.Range("D13").Value = 2028857
$ws.Range("E13").Value = 2467874
$ws.Range("F13").Value = 582654
$ws.Range("G13").Value = 1192787
$ws.Range("H13").Value = 1693506
$ws.Range("I13").Value = 1919214
$ws.Range("J13").Value = 899293
$ws.Range("K13").Value = 1730354
$ws.Range("L13").Value = 2537000
$ws.Range("M13").Value = 3604075

# ---------------------------------------------------------------------
# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
# ---------------------------------------------------------------------
$ws.Range("D14").Value = -115225
$ws.Range("E14").Value = -200101
$ws.Range("F14").Value = -61564
$ws.Range("G14").Value = -127635
$ws.Range("H14").Value = -204329
$ws.Range("I14").Value = -325003
$ws.Range("J14").Value = -83100
$ws.Range("K14").Value = -196015
$ws.Range("L14").Value = -321660
$ws.Range("M14").Value = -750939

# Row 15 (بهای استثنایی) is unchanged - all zero.

# ---------------------------------------------------------------------
# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 298830
$ws.Range("E16").Value = 281808
$ws.Range("F16").Value = 27772
$ws.Range("G16").Value = 118844
$ws.Range("H16").Value = 149030
$ws.Range("I16").Value = 129569
$ws.Range("J16").Value = 1138
$ws.Range("K16").Value = -4881
$ws.Range("L16").Value = 23725
$ws.Range("M16").Value = 165449

# ---------------------------------------------------------------------
# Row 17: سود (زیان) عملیاتی (Operating profit)
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 2212462
$ws.Range("E17").Value = 2549581
$ws.Range("F17").Value = 548862
$ws.Range("G17").Value = 1183996
$ws.Range("H17").Value = 1638207
$ws.Range("I17").Value = 1723780
$ws.Range("J17").Value = 817331
$ws.Range("K17").Value = 1529458
$ws.Range("L17").Value = 2239065
$ws.Range("M17").Value = 3018585

# ---------------------------------------------------------------------
# Row 18: هزینه های مالی (Finance expense)
# ---------------------------------------------------------------------
$ws.Range("D18").Value = -91674
$ws.Range("E18").Value = -137351
$ws.Range("F18").Value = -39255
$ws.Range("G18").Value = -96623
$ws.Range("H18").Value = -167434
$ws.Range("I18").Value = -276222
$ws.Range("J18").Value = -110158
$ws.Range("K18").Value = -230653
$ws.Range("L18").Value = -345067
$ws.Range("M18").Value = -462585

# ---------------------------------------------------------------------
# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
# ---------------------------------------------------------------------
$ws.Range("D19").Value = 17632
$ws.Range("E19").Value = 32820
$ws.Range("F19").Value = 17175
$ws.Range("G19").Value = 27627
$ws.Range("H19").Value = 35792
$ws.Range("I19").Value = 63918
$ws.Range("J19").Value = 1608
$ws.Range("K19").Value = 48330
$ws.Range("L19").Value = 56692
$ws.Range("M19").Value = 150319

# ---------------------------------------------------------------------
# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
# ---------------------------------------------------------------------
$ws.Range("D20").Value = 2138420
$ws.Range("E20").Value = 2445050
$ws.Range("F20").Value = 526782
$ws.Range("G20").Value = 1115000
$ws.Range("H20").Value = 1506565
$ws.Range("I20").Value = 1511476
$ws.Range("J20").Value = 708781
$ws.Range("K20").Value = 1347135
$ws.Range("L20").Value = 1950690
$ws.Range("M20").Value = 2706319

# ---------------------------------------------------------------------
# Row 21: مالیات (Tax)
# ---------------------------------------------------------------------
$ws.Range("D21").Value = -303999
$ws.Range("E21").Value = -287535
$ws.Range("F21").Value = -75077
$ws.Range("G21").Value = -209590
$ws.Range("H21").Value = -291004
$ws.Range("I21").Value = -274756
$ws.Range("J21").Value = -124253
$ws.Range("K21").Value = -94824
$ws.Range("L21").Value = -218219
$ws.Range("M21").Value = -122449

# ---------------------------------------------------------------------
# Row 22: سود (زیان) خالص عملیات در حال تداوم
# ---------------------------------------------------------------------
$ws.Range("D22").Value = 1834421
$ws.Range("E22").Value = 2157515
$ws.Range("F22").Value = 451705
$ws.Range("G22").Value = 905410
$ws.Range("H22").Value = 1215561
$ws.Range("I22").Value = 1236720
$ws.Range("J22").Value = 584528
$ws.Range("K22").Value = 1252311
$ws.Range("L22").Value = 1732471
$ws.Range("M22").Value = 2583870

# Row 23 (عملیات متوقف شده) is unchanged - all zero.

# ---------------------------------------------------------------------
# Row 24: سود (زیان) خالص (Net income)
# ---------------------------------------------------------------------
$ws.Range("D24").Value = 1834421
$ws.Range("E24").Value = 2157515
$ws.Range("F24").Value = 451705
$ws.Range("G24").Value = 905410
$ws.Range("H24").Value = 1215561
$ws.Range("I24").Value = 1236720
$ws.Range("J24").Value = 584528
$ws.Range("K24").Value = 1252311
$ws.Range("L24").Value = 1732471
$ws.Range("M24").Value = 2583870

# ---------------------------------------------------------------------
# Row 25: سود هر سهم پس از کسر مالیات (EPS)
# ---------------------------------------------------------------------
$ws.Range("D25").Value = 2038
$ws.Range("E25").Value = 2397
$ws.Range("F25").Value = 502
$ws.Range("G25").Value = 1006
$ws.Range("H25").Value = 1351
$ws.Range("I25").Value = 1374
$ws.Range("J25").Value = 649
$ws.Range("K25").Value = 1391
$ws.Range("L25").Value = 912
$ws.Range("M25").Value = 1360

# ---------------------------------------------------------------------
# Row 26: سرمایه (Capital) - only L26 changes (900000 -> 1900000)
# ---------------------------------------------------------------------
$ws.Range("L26").Value = 1900000

# ---------------------------------------------------------------------
# Row 27: سود هر سهم بر اساس آخرین سرمایه
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 965
$ws.Range("E27").Value = 1136
$ws.Range("F27").Value = 238
$ws.Range("G27").Value = 477
$ws.Range("H27").Value = 640
$ws.Range("I27").Value = 651
$ws.Range("J27").Value = 308
$ws.Range("K27").Value = 659
$ws.Range("L27").Value = 912
$ws.Range("M27").Value = 1360
